# Re-order / correct the observation rows so that each row's taxon data
# (Id, TaxonId, species names/author, coordinates, public comment) lines
# up with the right record. Net effect is a cyclic rotation of data among
# rows 2/3/4, a swap between rows 13/14, and a cyclic rotation among rows
# 17/18/19 (see commit diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData($Row, $A, $B, $E, $F, $G, $H, $Q, $R, $AC) {
    $ws.Range("A$Row").Value = $A
    $ws.Range("B$Row").Value = $B
    $ws.Range("E$Row").Value = $E
    $ws.Range("F$Row").Value = $F
    $ws.Range("G$Row").Value = $G

    if ($H) {
        $ws.Range("H$Row").Value = $H
    } else {
        $ws.Range("H$Row").ClearContents()
    }

    $ws.Range("Q$Row").Value = $Q
    $ws.Range("R$Row").Value = $R

    if ($AC) {
        $ws.Range("AC$Row").Value = $AC
    } else {
        $ws.Range("AC$Row").ClearContents()
    }
}

# Row 2 <- old row 4's data
Set-RowData 2 131064784 91828 5432 "Granticka" "Porodaedalea chrysoloma s.lat." $null 442100 7039221 $null

# Row 3 <- old row 2's data
Set-RowData 3 131064775 57884 100109 "Tretåig hackspett" "Picoides tridactylus" "(Linnaeus, 1758)" 442085 7039138 "Ringhack"

# Row 4 <- old row 3's data
Set-RowData 4 131064773 57884 100109 "Tretåig hackspett" "Picoides tridactylus" "(Linnaeus, 1758)" 442108 7039138 "Ringhack äldre"

# Row 13 <- old row 14's data
Set-RowData 13 131064779 91804 1108 "Harticka" "Pelloporus leporinus" "(Fr.) Krieglst." 442245 7039149 $null

# Row 14 <- old row 13's data
Set-RowData 14 131064763 57884 100109 "Tretåig hackspett" "Picoides tridactylus" "(Linnaeus, 1758)" 442230 7039147 "Ringhack äldre"

# Row 17 <- old row 18's data
Set-RowData 17 131064781 91804 1108 "Harticka" "Pelloporus leporinus" "(Fr.) Krieglst." 442200 7039150 $null

# Row 18 <- old row 19's data
Set-RowData 18 131064780 91804 1108 "Harticka" "Pelloporus leporinus" "(Fr.) Krieglst." 442259 7039181 $null

# Row 19 <- old row 17's data
Set-RowData 19 131064772 57884 100109 "Tretåig hackspett" "Picoides tridactylus" "(Linnaeus, 1758)" 442099 7039220 "Bohål ca 3m upp i grantickerötad granhögstubbe Även ett påbörjat på 2m"

Write-Output "Row rotation applied."
